# Re-order the header columns in row 2 (A2:L2) of the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newHeaders = @(
    "ProjectID",
    "BusinessKey",
    "OutcomeBusinessKey",
    "ProgrammeBusinessKey",
    "SectorBusinessKey",
    "SubSectorBusinessKey",
    "Code",
    "LongName",
    "ProjectParentID",
    "ProjectSiteName",
    "ShortName",
    "TextDescription"
)

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(2, $col).Value = $newHeaders[$i]
}
